# TC-61 and TC-62 added: insert two new worksheets ("tc061", "tc062") right
# before the "tc044" sheet, carrying the TR-notification test data, and
# adjust the view/selection state on the surrounding tabs (tc059, tc060)
# to match the post-edit session.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the two new sheets in the correct tab order --------------
# Worksheets.Add(Before) inserts immediately before the given sheet, so to
# end up with  ... tc060, tc061, tc062, tc044 ...  we add tc061 before
# tc044 first, then add tc062 right after tc061 (i.e. before tc044 too).
$sheetTc044 = $wb.Worksheets.Item("tc044")

$tc061 = $wb.Worksheets.Add($sheetTc044)
$tc061.Name = "tc061"

$tc062 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tc061)
$tc062.Name = "tc062"

# --- 2. Populate tc061 (8 columns: A:H) ----------------------------------
$tc061.Range("A1").Value = "Releasename"
$tc061.Range("B1").Value = "Cyclename"
$tc061.Range("C1").Value = "SubCyclename"
$tc061.Range("D1").Value = "Suitename"
$tc061.Range("E1").Value = "Epic"
$tc061.Range("F1").Value = "Feature"
$tc061.Range("G1").Value = "rq"
$tc061.Range("H1").Value = "tc"

$tc061.Range("A2").Value = "Release TR Notification 12-01-2026"
$tc061.Range("B2").Value = "TestCycle TR Notification 12-01-2026"
$tc061.Range("C2").Value = "Sub TestCycle TR Notification 12-01-2026"
$tc061.Range("D2").Value = "TestSuite TR Notification 12-01-2026"
$tc061.Range("E2").Value = "Epic Mohit"
$tc061.Range("F2").Value = "Mohit Feature"
$tc061.Range("G2").Value = "RQ-489"
$tc061.Range("H2").Value = "TC-427"

# --- 3. Populate tc062 (9 columns: A:I, adds a Status column) -----------
$tc062.Range("A1").Value = "Releasename"
$tc062.Range("B1").Value = "Cyclename"
$tc062.Range("C1").Value = "SubCyclename"
$tc062.Range("D1").Value = "Suitename"
$tc062.Range("E1").Value = "Epic"
$tc062.Range("F1").Value = "Feature"
$tc062.Range("G1").Value = "rq"
$tc062.Range("H1").Value = "tc"
$tc062.Range("I1").Value = "Status"

$tc062.Range("A2").Value = "Release TR Notification 12-01-2026"
$tc062.Range("B2").Value = "TestCycle TR Notification 12-01-2026"
$tc062.Range("C2").Value = "Sub TestCycle TR Notification 12-01-2026"
$tc062.Range("D2").Value = "TestSuite TR Notification 12-01-2026"
$tc062.Range("E2").Value = "Epic Mohit"
$tc062.Range("F2").Value = "Mohit Feature"
$tc062.Range("G2").Value = "RQ-489"
$tc062.Range("H2").Value = "TC-427"
$tc062.Range("I2").Value = "Passed"

# --- 4. View/selection tweaks on the touched tabs ------------------------
# tc059: selection grows from A1:A2 to the full A1:H2 data range.
$ws059 = $wb.Worksheets.Item("tc059")
$ws059.Activate()
$ws059.Range("A1:H2").Select()

# tc060: selection becomes the full A1:I2 data range (was a single cell).
$ws060 = $wb.Worksheets.Item("tc060")
$ws060.Activate()
$ws060.Range("A1:I2").Select()

# tc061: select the full range, landing the active cell on H2.
$tc061.Activate()
$tc061.Range("A1:H2").Select()
$tc061.Range("H2").Activate()

# tc062: final sheet worked on, so it ends up the active/selected tab.
$tc062.Activate()
$tc062.Range("E11").Select()

Write-Output "tc061 and tc062 added"
